$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the date/time value in A2 (existing row)
$ws.Range("A2").Value = 45956.56062929398

# Add a new expense row (row 3)
$ws.Range("A3").Value = 45957.86045812783
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B3").Value = "gasto"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = "mercado"
